$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.569.74"
$ws.Range("E2").Value = "  +1.75%  "

$ws.Range("D3").Value = "1.641.50"
$ws.Range("E3").Value = "  +2.58%  "

$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "307.83"
$ws.Range("E5").Value = "  +1.71%  "

$ws.Range("D6").Value = "0.9990"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").Value = "0.3762"
$ws.Range("E7").Value = "  -0.32%  "

$ws.Range("D8").Value = "52.75"
$ws.Range("E8").Value = "  +1.47%  "

$ws.Range("D9").Value = "0.3677"
$ws.Range("E9").Value = "  +1.99%  "

$ws.Range("D10").Value = "1.281"

$ws.Range("D11").Value = "0.08203"
$ws.Range("E11").Value = "  +1.09%  "

$ws.Range("D12").Value = "0.9981"
$ws.Range("E12").Value = "  +0.05%  "

$ws.Range("D13").Value = "23.09"
$ws.Range("E13").Value = "  +2.08%  "

$ws.Range("D14").Value = "6.687"
$ws.Range("E14").Value = "  +1.81%  "

$ws.Range("D15").Value = "0.00001286"
$ws.Range("E15").Value = "  +2.87%  "

$ws.Range("D16").Value = "7.462"
$ws.Range("E16").Value = "  +0.96%  "

$ws.Range("D17").Value = "1.639.38"
$ws.Range("E17").Value = "  +2.32%  "

$ws.Range("D18").Value = "95.10"
$ws.Range("E18").Value = "  +1.78%  "

$ws.Range("D19").Value = "0.06928"
$ws.Range("E19").Value = "  +1.26%  "

$ws.Range("D20").Value = "18.36"
$ws.Range("E20").Value = "  +1.98%  "

$ws.Range("D21").Value = "6.587"

$ws.Range("D22").Value = "0.9983"
$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("D23").Value = "23.603.08"
$ws.Range("E23").Value = "  +1.91%  "

$ws.Range("D24").Value = "12.92"
$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("D25").Value = "3.097"
$ws.Range("E25").Value = "  +4.06%  "

$ws.Range("D26").Value = "2.417"
$ws.Range("E26").Value = "  +1.39%  "

$ws.Range("D27").Value = "21.37"
$ws.Range("E27").Value = "  +1.09%  "

$ws.Range("D28").Value = "151.63"
$ws.Range("E28").Value = "  +1.44%  "

$ws.Range("D29").Value = "5.342"
$ws.Range("E29").Value = "  +2.34%  "

$ws.Range("D30").Value = "136.37"
$ws.Range("E30").Value = "  +1.90%  "

$ws.Range("D31").Value = "2.387"
$ws.Range("E31").Value = "  -0.33%  "

$ws.Range("D32").Value = "1.824.54"
$ws.Range("E32").Value = "  +2.73%  "

$ws.Range("D33").Value = "6.893"
$ws.Range("E33").Value = "  +1.19%  "

$ws.Range("D34").Value = "0.9827"
$ws.Range("E34").Value = "  +0.28%  "

$ws.Range("D35").Value = "0.02866"
$ws.Range("E35").Value = "  +5.62%  "

$ws.Range("D36").Value = "10.46"
$ws.Range("E36").Value = "  +1.59%  "

$ws.Range("D37").Value = "0.07399"
$ws.Range("E37").Value = "  -2.34%  "

$ws.Range("D38").Value = "0.2559"
$ws.Range("E38").Value = "  +2.27%  "

$ws.Range("D39").Value = "6.219"
$ws.Range("E39").Value = "  +1.29%  "

$ws.Range("D40").Value = "0.08905"

$ws.Range("D41").Value = "1.385"
$ws.Range("E41").Value = "  +1.72%  "

$ws.Range("D42").Value = "0.7154"
$ws.Range("E42").Value = "  +0.66%  "

$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "12.61"
$ws.Range("E43").Value = "  +1.98%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "16.34"
$ws.Range("E44").Value = "  +5.11%  "

$ws.Range("D45").Value = "0.6587"
$ws.Range("E45").Value = "  +1.17%  "

$ws.Range("D46").Value = "2.361"
$ws.Range("E46").Value = "  +2.74%  "

$ws.Range("D47").Value = "4.044"
$ws.Range("E47").Value = "  +0.90%  "

$ws.Range("D48").Value = "0.9980"
$ws.Range("E48").Value = "  +0.10%  "

$ws.Range("D49").Value = "0.08000"
$ws.Range("E49").Value = "  +0.74%  "

$ws.Range("D50").Value = "130.19"
$ws.Range("E50").Value = "  -1.50%  "

$ws.Range("E51").Value = "  +0.73%  "
